# Update "countries & provincias Spain" data
# - Refresh the "datos actualizados" timestamp
# - Paraguay overtakes Maldivas in the ranking (rows 108/109 swap position)
# - Islas Malvinas / Groenlandia tie-break swap (rows 209/210, identical values)
# - Updated case counts for several countries (rows 4, 6, 8, 27, 48, 54, 70, 94, 108, 109, 154)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp header (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Julio de 2020 a las 23:01"

# --- Row 4: Estados Unidos ---
$ws.Cells.Item(4, 2).Value = 3142533
$ws.Cells.Item(4, 3).Value = 45449
$ws.Cells.Item(4, 4).Value = 1379450
$ws.Cells.Item(4, 5).Value = 1628420
$ws.Cells.Item(4, 7).Value = 691
$ws.Cells.Item(4, 8).Value = 134663

# --- Row 6: India ---
$ws.Cells.Item(6, 2).Value = 769052
$ws.Cells.Item(6, 3).Value = 25571
$ws.Cells.Item(6, 4).Value = 476554
$ws.Cells.Item(6, 5).Value = 271354

# --- Row 8: Peru ---
$ws.Cells.Item(8, 2).Value = 312911
$ws.Cells.Item(8, 3).Value = 3633
$ws.Cells.Item(8, 4).Value = 204748
$ws.Cells.Item(8, 5).Value = 97030
$ws.Cells.Item(8, 7).Value = 181
$ws.Cells.Item(8, 8).Value = 11133

# --- Row 27: Egipto ---
$ws.Cells.Item(27, 2).Value = 78304
$ws.Cells.Item(27, 3).Value = 1025
$ws.Cells.Item(27, 4).Value = 22241
$ws.Cells.Item(27, 5).Value = 52499
$ws.Cells.Item(27, 7).Value = 75
$ws.Cells.Item(27, 8).Value = 3564

# --- Row 48: Israel ---
$ws.Cells.Item(48, 2).Value = 33557
$ws.Cells.Item(48, 3).Value = 1335
$ws.Cells.Item(48, 4).Value = 18338
$ws.Cells.Item(48, 5).Value = 14875

# --- Row 54: Irlanda ---
$ws.Cells.Item(54, 2).Value = 25542
$ws.Cells.Item(54, 3).Value = 11
$ws.Cells.Item(54, 5).Value = 440
$ws.Cells.Item(54, 8).Value = 1738

# --- Row 70: Costa de Marfil ---
$ws.Cells.Item(70, 2).Value = 11504
$ws.Cells.Item(70, 3).Value = 310
$ws.Cells.Item(70, 4).Value = 5571
$ws.Cells.Item(70, 5).Value = 5855
$ws.Cells.Item(70, 7).Value = 2
$ws.Cells.Item(70, 8).Value = 78

# --- Row 94: Estado de Palestina ---
$ws.Cells.Item(94, 5).Value = 4515
$ws.Cells.Item(94, 7).Value = 2
$ws.Cells.Item(94, 8).Value = 20

# --- Rows 108/109: Paraguay overtakes Maldivas ---
# Row 108 becomes Paraguay with updated figures
$ws.Cells.Item(108, 1).Value = "Paraguay"
$ws.Cells.Item(108, 2).Value = 2554
$ws.Cells.Item(108, 3).Value = 52
$ws.Cells.Item(108, 4).Value = 1212
$ws.Cells.Item(108, 5).Value = 1322
$ws.Cells.Item(108, 7).Value = 0
$ws.Cells.Item(108, 8).Value = 20

# Row 109 becomes Maldivas, carrying the figures that used to be in row 108
$ws.Cells.Item(109, 1).Value = "Maldivas"
$ws.Cells.Item(109, 2).Value = 2517
$ws.Cells.Item(109, 3).Value = 16
$ws.Cells.Item(109, 4).Value = 2180
$ws.Cells.Item(109, 5).Value = 324
$ws.Cells.Item(109, 7).Value = 1
$ws.Cells.Item(109, 8).Value = 13

# --- Rows 209/210: Islas Malvinas / Groenlandia swap (values identical, order ties) ---
$ws.Cells.Item(209, 1).Value = "Islas Malvinas"
$ws.Cells.Item(210, 1).Value = "Groenlandia"

# --- Row 154: Surinam ---
$ws.Cells.Item(154, 2).Value = 665
$ws.Cells.Item(154, 3).Value = 31
$ws.Cells.Item(154, 4).Value = 434
$ws.Cells.Item(154, 7).Value = 2
$ws.Cells.Item(154, 8).Value = 17
